$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.679.28"
$ws.Range("D3").Value = "2.968.72"
$ws.Range("E3").Value = "  -6.22%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'541.54"
$ws.Range("E5").Value = "  -5.33%  "
$ws.Range("D6").Value = "'152.87"
$ws.Range("E6").Value = "  -6.99%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -3.93%  "
$ws.Range("D9").Value = "2.977.17"
$ws.Range("E9").Value = "  -5.86%  "
$ws.Range("D10").Value = "'0.113"
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("D11").Value = "'6.10"
$ws.Range("E11").Value = "  -8.08%  "
$ws.Range("D12").Value = "'0.365"
$ws.Range("E12").Value = "  -5.30%  "
$ws.Range("D13").Value = "3.488.66"
$ws.Range("E13").Value = "  -6.11%  "
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "61.710.03"
$ws.Range("E15").Value = "  -4.24%  "
$ws.Range("D16").Value = "'23.64"
$ws.Range("E16").Value = "  -6.77%  "
$ws.Range("D17").Value = "2.973.13"
$ws.Range("E17").Value = "  -6.04%  "
$ws.Range("E18").Value = "  -6.34%  "
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'11.95"
$ws.Range("E20").Value = "  -6.53%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'379.37"
$ws.Range("E21").Value = "  -7.38%  "
$ws.Range("E22").Value = "  -6.35%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'65.06"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("D26").Value = "3.087.86"
$ws.Range("E26").Value = "  -6.52%  "
$ws.Range("E27").Value = "  -5.77%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").Value = "0.0₃0931"
$ws.Range("E29").Value = "  -9.63%  "
$ws.Range("E30").Value = "  -7.58%  "
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  -5.47%  "
$ws.Range("D33").Value = "'20.46"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "'158.78"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  -5.59%  "
$ws.Range("D36").Value = "'5.96"
$ws.Range("E36").Value = "  -6.40%  "
$ws.Range("E37").Value = "  -6.33%  "
$ws.Range("E38").Value = "  -5.71%  "
$ws.Range("D39").Value = "'1.56"
$ws.Range("E39").Value = "  -8.48%  "
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").Value = "2.411.79"
$ws.Range("E41").Value = "  -10.29%  "
$ws.Range("D42").Value = "'37.06"
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").Value = "'22.14"
$ws.Range("E43").Value = "  -7.99%  "
$ws.Range("E44").Value = "  -4.49%  "
$ws.Range("D45").Value = "'0.0589"
$ws.Range("E45").Value = "  -5.13%  "
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("E47").Value = "  -5.43%  "
$ws.Range("D48").Value = "'5.01"
$ws.Range("E48").Value = "  -8.26%  "
$ws.Range("D49").Value = "'0.0945"
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("D50").Value = "'19.61"
$ws.Range("E50").Value = "  -8.94%  "
$ws.Range("E51").Value = "  +0.17%  "
